$wb = $excel.ActiveWorkbook

# Update the zh-cn sheet: Correspond Handoff Datetime (D2) and
# Correspond Handback DateTime (G2) for the first data row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-23 07:39:29"
$wsZhCn.Range("G2").Value = "2016-02-23 07:40:18"

# Update the de-de sheet: Correspond Handoff Datetime (D2) and
# Correspond Handback DateTime (G2) for the first data row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-23 07:39:43"
$wsDeDe.Range("G2").Value = "2016-02-23 07:40:44"
